$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D range to Text format so numeric-looking strings
# (e.g. "242.65") are stored as literal text, matching the source data,
# which uses text values throughout (t="inlineStr") rather than numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.084.58"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.836.62"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "242.65"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "0.6261"
$ws.Range("E6").Value = "  -5.53%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "0.07554"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").Value = "22.64"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").Value = "0.07751"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.840.68"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "4.954"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "82.78"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "0.000009926"
$ws.Range("E16").Value = "  +13.58%  "
$ws.Range("D17").Value = "6.031"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").Value = "29.124.40"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "225.76"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "7.182"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "158.77"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "8.467"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").Value = "0.1366"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "17.91"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "1.490"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").Value = "4.078"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").Value = "4.034"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "1.201"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").Value = "0.05204"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").Value = "1.857"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").Value = "0.7387"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "1.250.86"
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("D38").Value = "2.765"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").Value = "6.370"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").Value = "0.8928"
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "101.65"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "1.980.08"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "0.00000000125"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "64.06"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "0.5121"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "0.4019"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "8.844"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").Value = "1.646"
$ws.Range("E50").Value = "  -6.31%  "
$ws.Range("D51").Value = "0.05761"
$ws.Range("E51").Value = "  -1.81%  "

# Restore the original (default/Normal) cell style on column D so we
# do not leave a lingering text-format style applied to the cells.
$dRange.Style = "Normal"
